$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "30.351.09"
$ws.Range("E2").Value = "  -2.54%  "
Set-TextValue $ws.Range("D3") "1.943.65"
$ws.Range("E3").Value = "  -2.32%  "
Set-TextValue $ws.Range("D4") "1.000"
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue $ws.Range("D5") "252.26"
Set-TextValue $ws.Range("D6") "0.7246"
$ws.Range("E6").Value = "  -8.06%  "
Set-TextValue $ws.Range("D7") "1.000"
$ws.Range("E7").Value = "  -0.02%  "
Set-TextValue $ws.Range("D8") "0.3350"
$ws.Range("E8").Value = "  -4.11%  "
Set-TextValue $ws.Range("D9") "28.87"
$ws.Range("E9").Value = "  +3.28%  "
Set-TextValue $ws.Range("D10") "0.07457"
$ws.Range("E10").Value = "  +6.82%  "
Set-TextValue $ws.Range("D11") "0.8221"
$ws.Range("E11").Value = "  -2.46%  "
$ws.Range("E12").Value = "  -0.59%  "
Set-TextValue $ws.Range("D13") "1.940.75"
$ws.Range("E13").Value = "  -2.48%  "
Set-TextValue $ws.Range("D14") "5.506"
$ws.Range("E14").Value = "  -1.52%  "
Set-TextValue $ws.Range("D15") "95.49"
$ws.Range("E15").Value = "  -4.68%  "
Set-TextValue $ws.Range("D16") "14.97"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D17") "0.000008410"
$ws.Range("E17").Value = "  +6.56%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D18") "30.361.21"
$ws.Range("E18").Value = "  -2.52%  "
Set-TextValue $ws.Range("D19") "254.43"
$ws.Range("E19").Value = "  -7.00%  "
Set-TextValue $ws.Range("D20") "5.912"
$ws.Range("E20").Value = "  +0.69%  "
Set-TextValue $ws.Range("D21") "2.196.61"
$ws.Range("E21").Value = "  -2.55%  "
$ws.Range("E22").Value = "  -0.04%  "
Set-TextValue $ws.Range("D23") "0.9983"
$ws.Range("E23").Value = "  -0.27%  "
Set-TextValue $ws.Range("D24") "7.009"
$ws.Range("E24").Value = "  -0.65%  "
Set-TextValue $ws.Range("D25") "9.902"
$ws.Range("E25").Value = "  -1.38%  "
Set-TextValue $ws.Range("D26") "162.33"
$ws.Range("E26").Value = "  -1.42%  "
Set-TextValue $ws.Range("D27") "2.431"
$ws.Range("E27").Value = "  +4.71%  "
Set-TextValue $ws.Range("D28") "19.44"
$ws.Range("E28").Value = "  -2.38%  "
$ws.Range("E29").Value = "  -11.78%  "
$ws.Range("E30").Value = "  -1.42%  "
Set-TextValue $ws.Range("D31") "1.347"
$ws.Range("E31").Value = "  -1.03%  "
Set-TextValue $ws.Range("D32") "4.476"
Set-TextValue $ws.Range("D33") "4.272"
$ws.Range("E33").Value = "  -3.11%  "
Set-TextValue $ws.Range("D34") "0.05320"
$ws.Range("E34").Value = "  +1.99%  "
Set-TextValue $ws.Range("D35") "1.313"
$ws.Range("E35").Value = "  +7.19%  "
Set-TextValue $ws.Range("D36") "0.7630"
$ws.Range("E36").Value = "  -1.98%  "
Set-TextValue $ws.Range("D37") "2.749"
$ws.Range("E37").Value = "  -0.35%  "
Set-TextValue $ws.Range("D38") "0.01999"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("E39").Value = "  -1.56%  "
Set-TextValue $ws.Range("D40") "81.52"
$ws.Range("E40").Value = "  +2.89%  "
Set-TextValue $ws.Range("D41") "6.618"
$ws.Range("E41").Value = "  -0.02%  "
Set-TextValue $ws.Range("D42") "0.4588"
$ws.Range("E42").Value = "  -1.47%  "
Set-TextValue $ws.Range("D43") "2.050"
$ws.Range("E43").Value = "  -3.19%  "
Set-TextValue $ws.Range("D44") "0.8450"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("E45").Value = "  +0.01%  "
Set-TextValue $ws.Range("D46") "103.19"
$ws.Range("E46").Value = "  -1.76%  "
Set-TextValue $ws.Range("D47") "9.831"
$ws.Range("E47").Value = "  -0.85%  "
Set-TextValue $ws.Range("D48") "7.536"
$ws.Range("E48").Value = "  -1.73%  "
Set-TextValue $ws.Range("D49") "37.11"
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("E50").Value = "  -1.72%  "
Set-TextValue $ws.Range("D51") "1.523"
$ws.Range("E51").Value = "  -0.37%  "
